$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell C10: value changes from 18 to 1
$ws.Range("C10").Value = 1
